$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The derived-data export no longer carries the computed "Sphericity" (G)
# and "SHE" (H) columns -- drop them outright. This shifts the trailing
# "Cortex covered" column left so it becomes the new column G.
$ws.Range("G:H").EntireColumn.Delete()

# Keep the sheet's hidden AutoFilter range (_FilterDatabase) in sync with
# the new, narrower used range instead of leaving it pointing past the
# data that no longer exists.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase") {
        $n.RefersTo = "=" + $ws.Name + "!" + $ws.UsedRange.Address()
    }
}
